# "Generate Report for Archive"
# - Status text "Ready for handoff" -> "In Translation" (all sheets)
# - Narrow the "Status"/"zh-cn"/"de-de" status columns (report regenerated
#   with a shorter status string, so the column auto-sized narrower)

$wb = $excel.ActiveWorkbook

# 1) Update the status text everywhere it appears (Overview!E:F, zh-cn!C, de-de!C)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2) Shrink the status columns to their new (narrower) width.
#    Target stored width ~= 13.41 characters; this engine quantizes
#    ColumnWidth to 1/6-character steps, so 12.42 is the closest input
#    that lands in that bucket.
$newColumnWidth = 12.42

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
